$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update the "exit_size" values in column E (rows 3-5) ---
# Order matters: it rebuilds the shared-string table positions to match
# the target layout (indices 50/51/52 => 2.6, 3.42 / 3.23, 4.79 / 6.41, 9.55)
$ws.Range("E5").Value = "2.6, 3.42"
$ws.Range("E4").Value = "3.23, 4.79"
$ws.Range("E3").Value = "6.41, 9.55"

# --- Correct F4 (wallthick) value ---
$ws.Range("F4").Value = 0.075

# --- Add the new "used to be" comparison header block (G1, I1:K1) ---
$ws.Range("G1").Value = "used to be"
$ws.Range("I1").Value = "exit_size"
$ws.Range("J1").Value = "slits"
$ws.Range("K1").Value = "wallthick"

# --- Highlight the I2:K5 block with a solid red fill ---
$ws.Range("I2:K5").Interior.Color = 255

# --- Widen column J (10th column) to fit the new header text ---
$ws.Columns.Item(10).ColumnWidth = 14.9

# --- Update the last active selection to E3 ---
$ws.Range("E3").Select() | Out-Null
